$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header-cell formatting (bold, border, centered) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: column I (I0) and column J (IF)
$values = @{
    2  = @(1, 7)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 7)
    10 = @(1, 5)
    11 = @(6, 7)
    12 = @(8, 9)
    13 = @(1, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
